# Updates cryptos list with latest price/volume snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.749.98"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "'2.230.57"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'313.96"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "'98.10"
$ws.Range("E6").Value = "  -4.56%  "
$ws.Range("D7").Value = "'0.569"
$ws.Range("E7").Value = "  -3.01%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -6.84%  "
$ws.Range("D10").Value = "'36.14"
$ws.Range("E10").Value = "  -6.45%  "
$ws.Range("D11").Value = "'0.0824"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "'7.40"
$ws.Range("E12").Value = "  -5.58%  "
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "'2.566.15"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.839"
$ws.Range("E15").Value = "  -4.28%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.231.25"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "'14.11"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("D18").Value = "'43.620.83"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'12.91"
$ws.Range("E19").Value = "  -10.99%  "
$ws.Range("D20").Value = "'0.0₃0964"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  -4.92%  "
$ws.Range("D22").Value = "'65.20"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").Value = "'2.99"
$ws.Range("E23").Value = "  -7.17%  "
$ws.Range("D24").Value = "'233.58"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "'2.03"
$ws.Range("E25").Value = "  -7.37%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "'2.17"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("D29").Value = "'36.42"
$ws.Range("E29").Value = "  -7.26%  "
$ws.Range("D30").Value = "'5.95"
$ws.Range("E30").Value = "  -9.27%  "
$ws.Range("D31").Value = "'158.92"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'19.92"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "'0.0831"
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "'3.21"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.90"
$ws.Range("E36").Value = "  -5.29%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.109"
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").Value = "'15.78"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  -7.35%  "
$ws.Range("D41").Value = "'4.02"
$ws.Range("E41").Value = "  -11.95%  "
$ws.Range("D42").Value = "'0.0307"
$ws.Range("E42").Value = "  -5.85%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'1.716.03"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").Value = "'0.194"
$ws.Range("E45").Value = "  -6.56%  "
$ws.Range("D46").Value = "'81.12"
$ws.Range("E46").Value = "  -5.04%  "
$ws.Range("D47").Value = "'73.41"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("D48").Value = "'5.11"
$ws.Range("E48").Value = "  -5.38%  "
$ws.Range("D49").Value = "'102.31"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "'56.46"
$ws.Range("E51").Value = "  -5.45%  "
